# "feat : existance d'adhérent"
#
# Reorders / edits the four "Adhérent :" bullet items (numId 7) under the
# "Prêter un exemplaire" use case:
#   1. "Exemplaire disponible"                -> append " (OK)"
#   2. "Quota de prêt (... ) > 0 "             -> becomes "Adhérent existe vraiment (OK)"
#   3. "Adhérent existe vraiment"              -> becomes "Adhérent abonnée"
#   4. "Adhérent abonnée"                      -> becomes the old "Quota de prêt (...) > 0 " text
#
# i.e. the "Quota de prêt" bullet moves down two slots (after "Adhérent
# abonnée"), "Adhérent existe vraiment" moves up one slot and gets a
# "(OK)" marker, and "Exemplaire disponible" also gets a "(OK)" marker.

$d = $word.ActiveDocument

$quotaText = "Quota de prêt (nombre de livre qu" + [char]0x2019 + "un adhérent peut prêter simultanément par rapport à son profil) > 0 "

# Locate the 4 consecutive paragraphs by content so the script is resilient
# to the exact paragraph index.
$idxDisponible = -1
$idxQuota = -1
$idxExiste = -1
$idxAbonnee = -1

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($idxDisponible -eq -1 -and $t -like "Exemplaire disponible*") {
        $idxDisponible = $i
    } elseif ($idxQuota -eq -1 -and $t -like "Quota de*") {
        $idxQuota = $i
    } elseif ($idxExiste -eq -1 -and $t -like "Adhérent existe*") {
        $idxExiste = $i
    } elseif ($idxAbonnee -eq -1 -and $t -like "Adhérent abonn*") {
        $idxAbonnee = $i
    }
}

if ($idxDisponible -eq -1 -or $idxQuota -eq -1 -or $idxExiste -eq -1 -or $idxAbonnee -eq -1) {
    throw "Could not locate all four target paragraphs (Disponible=$idxDisponible Quota=$idxQuota Existe=$idxExiste Abonnee=$idxAbonnee)"
}

# 1. "Exemplaire disponible" -> append " (OK)"
$r = $d.Paragraphs.Item($idxDisponible).Range
$r.MoveEnd(1, -1) | Out-Null
$r.InsertAfter(" (OK)")

# 2. "Quota de prêt ..." paragraph -> "Adhérent existe vraiment (OK)"
$r = $d.Paragraphs.Item($idxQuota).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "Adhérent existe vraiment (OK)"

# 3. "Adhérent existe vraiment" paragraph -> "Adhérent abonnée"
$r = $d.Paragraphs.Item($idxExiste).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "Adhérent abonnée"

# 4. "Adhérent abonnée" paragraph -> the original "Quota de prêt (...) > 0 " text
$r = $d.Paragraphs.Item($idxAbonnee).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = $quotaText

Write-Host "Done."
